{"js": "// Apply \"resultados atualizados com dados corrigidos\": update the\n// comorbidades table header labels and the refreshed counts/percentages.\nconst replacements = [\n  [\"categoria\", \"Categoria\"],\n  [\"Global\", \"Quantidade\"],\n  [\"422\", \"365\"],\n  [\"407 (96.4)\", \"353 (96.7)\"],\n  [\"15 ( 3.6)\", \"12 ( 3.3)\"],\n  [\"409 (96.9)\", \"355 (97.3)\"],\n  [\"13 ( 3.1)\", \"10 ( 2.7)\"],\n  [\"352 (95.4)\", \"299 (95.8)\"],\n  [\"17 ( 4.6)\", \"13 ( 4.2)\"],\n  [\"376 (89.1)\", \"326 (89.3)\"],\n  [\"46 (10.9)\", \"39 (10.7)\"],\n  [\"418 (99.1)\", \"363 (99.5)\"],\n  [\"4 ( 0.9)\", \"2 ( 0.5)\"],\n  [\"195 (46.2)\", \"167 (45.8)\"],\n  [\"227 (53.8)\", \"198 (54.2)\"],\n  [\"327 (77.5)\", \"286 (78.4)\"],\n  [\"95 (22.5)\", \"79 (21.6)\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply \"resultados atualizados com dados corrigidos\": update the\n# comorbidades table header labels and the refreshed counts/percentages.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"categoria\", \"Categoria\"),\n  @(\"Global\", \"Quantidade\"),\n  @(\"422\", \"365\"),\n  @(\"407 (96.4)\", \"353 (96.7)\"),\n  @(\"15 ( 3.6)\", \"12 ( 3.3)\"),\n  @(\"409 (96.9)\", \"355 (97.3)\"),\n  @(\"13 ( 3.1)\", \"10 ( 2.7)\"),\n  @(\"352 (95.4)\", \"299 (95.8)\"),\n  @(\"17 ( 4.6)\", \"13 ( 4.2)\"),\n  @(\"376 (89.1)\", \"326 (89.3)\"),\n  @(\"46 (10.9)\", \"39 (10.7)\"),\n  @(\"418 (99.1)\", \"363 (99.5)\"),\n  @(\"4 ( 0.9)\", \"2 ( 0.5)\"),\n  @(\"195 (46.2)\", \"167 (45.8)\"),\n  @(\"227 (53.8)\", \"198 (54.2)\"),\n  @(\"327 (77.5)\", \"286 (78.4)\"),\n  @(\"95 (22.5)\", \"79 (21.6)\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  $found = $find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n  if (-not $found) {\n    throw \"Text not found: $($pair[0])\"\n  }\n}\n"}
